$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.623.06'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.843.58'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.72%  '
$ws.Range('E4').Value = '  -0.63%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.552'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('E8').Value = '  +5.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.293'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0712'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.59%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.099.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.12%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.20'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.837.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.655'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.625.56'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '254.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0804'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.996'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.65%  '
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.85%  '
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('E30').Value = '  +3.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.04%  '
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '504.52'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +878.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.97%  '
$ws.Range('E35').Value = '  +6.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.442.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('E37').Value = '  +4.83%  '
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('E39').Value = '  +3.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.988'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.12%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '83.23'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.13%  '
$ws.Range('E45').Value = '  +5.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.993.69'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.54%  '
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('E49').Value = '  -3.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '106.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.69%  '
$ws.Range('E51').Value = '  -0.23%  '
